$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the number format currently used by the "last row" date cell (A28)
# before we touch anything, and the format used by a regular (non-last) date
# cell (A27), so we can swap them around as the new row is appended.
$lastRowDateFormat = $ws.Range("A28").NumberFormat
$regularDateFormat = $ws.Range("A27").NumberFormat

# Row 28 is no longer the last row, so it reverts to the regular date format.
$ws.Range("A28").NumberFormat = $regularDateFormat

# Append the new day's data as row 29.
$ws.Range("A29").Value = 45978
$ws.Range("B29").Value = 63
$ws.Range("C29").Value = 71
$ws.Range("D29").Value = 73

# Row 29, now the last row, takes on the "last row" date format.
$ws.Range("A29").NumberFormat = $lastRowDateFormat
